# New crime data collected - update the weekly CompStat report
# (66th Precinct, week of 6/26/2023 through 7/2/2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  25" -> "Volume 30   Number  26"
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  26"

# "Report Covering the Week  6/19/2023  Through  6/25/2023"
#   -> "Report Covering the Week  6/26/2023  Through  7/2/2023"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Crime complaint data table (rows 16-30, columns C-N) ---------------
# Columns: C=WTD 2023, D=WTD 2022, E=WTD %Chg, F=28Day 2023, G=28Day 2022,
#          H=28Day %Chg, I=YTD 2023, J=YTD 2022, K=YTD %Chg, L=2Yr %Chg,
#          M=13Yr %Chg, N=30Yr %Chg

# Row 16 - Robbery
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 6
$ws.Cells.Item(16, 7).Value = 8
$ws.Cells.Item(16, 8).Value = -25
$ws.Cells.Item(16, 9).Value = 40
$ws.Cells.Item(16, 10).Value = 42
$ws.Cells.Item(16, 11).Value = -4.761904761904
$ws.Cells.Item(16, 12).Value = 33.333333333333
$ws.Cells.Item(16, 13).Value = -50.617283950617
$ws.Cells.Item(16, 14).Value = -89.189189189189

# Row 17 - Fel. Assault
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 21
$ws.Cells.Item(17, 7).Value = 17
$ws.Cells.Item(17, 8).Value = 23.529411764705
$ws.Cells.Item(17, 9).Value = 107
$ws.Cells.Item(17, 10).Value = 105
$ws.Cells.Item(17, 11).Value = 1.904761904761
$ws.Cells.Item(17, 12).Value = 62.121212121212
$ws.Cells.Item(17, 13).Value = 59.701492537313
$ws.Cells.Item(17, 14).Value = -31.847133757961

# Row 18 - Burglary
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 7).Value = 17
$ws.Cells.Item(18, 8).Value = -41.176470588235
$ws.Cells.Item(18, 9).Value = 60
$ws.Cells.Item(18, 10).Value = 84
$ws.Cells.Item(18, 11).Value = -28.571428571428
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -68.75
$ws.Cells.Item(18, 14).Value = -92.736077481840

# Row 19 - Gr. Larceny
$ws.Cells.Item(19, 3).Value = 13
$ws.Cells.Item(19, 4).Value = 21
$ws.Cells.Item(19, 5).Value = -38.095238095238
$ws.Cells.Item(19, 6).Value = 45
$ws.Cells.Item(19, 7).Value = 68
$ws.Cells.Item(19, 8).Value = -33.823529411764
$ws.Cells.Item(19, 9).Value = 277
$ws.Cells.Item(19, 10).Value = 297
$ws.Cells.Item(19, 11).Value = -6.734006734006
$ws.Cells.Item(19, 12).Value = 31.904761904761
$ws.Cells.Item(19, 13).Value = 40.609137055837
$ws.Cells.Item(19, 14).Value = -9.771986970684

# Row 20 - G.L.A.
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 4).Value = 5
$ws.Cells.Item(20, 5).Value = 20
$ws.Cells.Item(20, 6).Value = 17
$ws.Cells.Item(20, 7).Value = 10
$ws.Cells.Item(20, 8).Value = 70
$ws.Cells.Item(20, 9).Value = 77
$ws.Cells.Item(20, 10).Value = 48
$ws.Cells.Item(20, 11).Value = 60.416666666666
$ws.Cells.Item(20, 12).Value = 102.631578947368
$ws.Cells.Item(20, 13).Value = 6.944444444444
$ws.Cells.Item(20, 14).Value = -91.920251836306

# Row 21 - TOTAL
$ws.Cells.Item(21, 3).Value = 26
$ws.Cells.Item(21, 4).Value = 32
$ws.Cells.Item(21, 5).Value = -18.75
$ws.Cells.Item(21, 6).Value = 99
$ws.Cells.Item(21, 7).Value = 121
$ws.Cells.Item(21, 8).Value = -18.181818181818
$ws.Cells.Item(21, 9).Value = 572
$ws.Cells.Item(21, 10).Value = 588
$ws.Cells.Item(21, 11).Value = -2.721088435374
$ws.Cells.Item(21, 12).Value = 39.853300733496
$ws.Cells.Item(21, 13).Value = -6.382978723404
$ws.Cells.Item(21, 14).Value = -78.234398782344

# Row 24 - Petit Larceny
$ws.Cells.Item(24, 3).Value = 23
$ws.Cells.Item(24, 4).Value = 15
$ws.Cells.Item(24, 5).Value = 53.333333333333
$ws.Cells.Item(24, 7).Value = 92
$ws.Cells.Item(24, 8).Value = 4.347826086956
$ws.Cells.Item(24, 9).Value = 577
$ws.Cells.Item(24, 10).Value = 555
$ws.Cells.Item(24, 11).Value = 3.963963963963
$ws.Cells.Item(24, 12).Value = 65.804597701149
$ws.Cells.Item(24, 13).Value = 29.082774049217

# Row 25 - Misd. Assault
$ws.Cells.Item(25, 3).Value = 10
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = 150
$ws.Cells.Item(25, 6).Value = 35
$ws.Cells.Item(25, 7).Value = 21
$ws.Cells.Item(25, 8).Value = 66.666666666666
$ws.Cells.Item(25, 9).Value = 172
$ws.Cells.Item(25, 10).Value = 160
$ws.Cells.Item(25, 11).Value = 7.5
$ws.Cells.Item(25, 12).Value = 27.407407407407
$ws.Cells.Item(25, 13).Value = -24.561403508771

# Row 27 - Other Sex Crimes
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = 33.333333333333
$ws.Cells.Item(27, 9).Value = 33
$ws.Cells.Item(27, 10).Value = 38
$ws.Cells.Item(27, 11).Value = -13.157894736842
$ws.Cells.Item(27, 12).Value = 6.451612903225

# Row 28 - Shooting Vic.
$ws.Cells.Item(28, 14).Value = -80

# Row 29 - Shooting Inc.
$ws.Cells.Item(29, 14).Value = -80

# Row 30 - Hate Crimes (C30 and F30 go from "no data" text placeholder to
# real numeric values, so the number format needs to be (re)applied too)
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 3).NumberFormat = "#,##0"
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 6).NumberFormat = "#,##0"
$ws.Cells.Item(30, 8).Value = -50
$ws.Cells.Item(30, 9).Value = 4
$ws.Cells.Item(30, 11).Value = -55.555555555555
$ws.Cells.Item(30, 12).Value = -42.857142857142
